$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "paris"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "CRT"
$ws.Range("E5").Value = "RES"
$ws.Range("F5").Value = "42be9703-0e9b-4ce8-962d-60bf1f233ce8"
$ws.Range("G5").Value = "SJCPLLpaW_annotated.xlsx"
$ws.Range("H5").Value = "The results show that DeePa achieves speedups compared to PyTorch and TensorFlow with all of the tested minibatch sizes."
